$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 28 with data values
$ws.Range("A28").Value = 1
$ws.Range("B28").Value = "chr"
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = "chr"

# Update the view: scroll so row 8 is the top row, and move the active selection to A28
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("A28").Select()
